$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data fix: row 29 parent-1 name was shortened ---
$ws.Range("C29").Value = "אורלי"

# --- New row 33: a teacher entry added to the phonebook ---
$ws.Range("A33").Value = "המורה יפית"
$ws.Range("B33").Value = "קינן"
$ws.Range("D33").Value = 972528092008
$ws.Range("H33").Value = "teacher.jpg"
$ws.Range("J33").Value = "girl"

# --- Extend the gender list-validation down to the new row ---
$ws.Range("J2:J32").Validation.Delete() | Out-Null
$ws.Range("J2:J33").Validation.Add(3, 1, 1, "=`$M`$3:`$M`$4") | Out-Null

# --- Sheet view: scroll/selection moved while editing the new row ---
$ws.Range("J32:J33").Select() | Out-Null
